$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.227.75'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = '2.073.29'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9985'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5185'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4308'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08820'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.152'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").Value = '2.071.65'
$ws.Range("E13").Value = '  +3.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.646'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.657'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.77%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001120'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06597'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9994'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.205'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '30.279.46'
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.275'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("D26").Value = '2.313.80'
$ws.Range("E26").Value = '  +3.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.526'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.29%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '130.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.187'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.077'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.551'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +14.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.833'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02558'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.593'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.11%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06595'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.394'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2227'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6656'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.238'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9982'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6301'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.187'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.177'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.10%  '
